$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.643.20'
$ws.Range("E2").Value = '  -3.04%  '

$ws.Range("D3").Value = '2.614.07'
$ws.Range("E3").Value = '  -1.34%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.12'
$ws.Range("E5").Value = '  -3.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.74'
$ws.Range("E6").Value = '  -1.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.630'
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("E9").Value = '  -5.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.78'
$ws.Range("E10").Value = '  -0.72%  '

$ws.Range("E11").Value = '  -3.03%  '

$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.16'
$ws.Range("E13").Value = '  -1.65%  '

$ws.Range("D14").Value = '3.085.21'
$ws.Range("E14").Value = '  -1.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000183'
$ws.Range("E15").Value = '  -7.47%  '

$ws.Range("D16").Value = '63.566.48'
$ws.Range("E16").Value = '  -2.91%  '

$ws.Range("D17").Value = '2.637.96'
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.04'
$ws.Range("E18").Value = '  -4.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.62'
$ws.Range("E19").Value = '  -2.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.45'
$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.25'
$ws.Range("E21").Value = '  -1.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.98'
$ws.Range("E23").Value = '  -2.86%  '

$ws.Range("E24").Value = '  +2.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000108'
$ws.Range("E25").Value = '  -4.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.23'
$ws.Range("E26").Value = '  -3.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '578.36'
$ws.Range("E27").Value = '  +9.20%  '

$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.161'
$ws.Range("E29").Value = '  -2.08%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.70'
$ws.Range("E33").Value = '  -3.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.47'
$ws.Range("E34").Value = '  +0.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.31'
$ws.Range("E35").Value = '  -1.65%  '

$ws.Range("E36").Value = '  -2.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.85'
$ws.Range("E37").Value = '  -2.49%  '

$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '153.46'
$ws.Range("E39").Value = '  -1.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.86'
$ws.Range("E40").Value = '  -3.37%  '

$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.32'
$ws.Range("E42").Value = '  -2.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '156.10'
$ws.Range("E43").Value = '  -3.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.37'
$ws.Range("E44").Value = '  +4.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.95'
$ws.Range("E45").Value = '  -2.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.87'
$ws.Range("E46").Value = '  +1.46%  '

$ws.Range("E47").Value = '  -2.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.628'
$ws.Range("E48").Value = '  -1.10%  '

$ws.Range("E49").Value = '  +1.99%  '

$ws.Range("E50").Value = '  -1.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.90'
$ws.Range("E51").Value = '  -4.15%  '
